# Box Plot Updates, Color Updates Main Figures
#
# Slide 1 contains a single top-level group shape (Shapes.Item(2)) that
# holds the chart pieces; the label textboxes tx9..tx18 live inside it as
# GroupItems 7..16. This nudges their positions to match the refreshed
# layout (some move only slightly to re-align with the redrawn pie/box
# plot wedges, two (tx15/tx16, the "Heterotroph" callout) move further to
# track the relocated wedge).
#
# NOTE: PowerPoint's Shape.Left/.Top are COM `Single` (32-bit float)
# properties expressed in points (1 pt = 12700 EMU). A naive
# `emu / 12700.0` assignment can truncate to one EMU below the intended
# value once it round-trips through the 32-bit float, so the literals
# below are the nearest points value that reproduces the exact target
# EMU after that conversion.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$g = $s.Shapes.Item(2)

# tx9 "Plastid Parasite" : (4505257,2460508) -> (4505183,2460418)
$sh = $g.GroupItems.Item(7)
$sh.Left = 354.73883056640625
$sh.Top  = 193.73370361328125

# tx10 " 2.459 %" : (4860644,2766892) -> (4860569,2766802)
$sh = $g.GroupItems.Item(8)
$sh.Left = 382.72198486328125
$sh.Top  = 217.85842895507812

# tx11 "Non-Plastid Parasite" : (5246207,3068257) -> (5246240,3068225)
$sh = $g.GroupItems.Item(9)
$sh.Left = 413.08978271484375
$sh.Top  = 241.592529296875

# tx12 " 7.377 %" : (5836559,3374641) -> (5836592,3374608)
$sh = $g.GroupItems.Item(10)
$sh.Left = 459.5741882324219
$sh.Top  = 265.7171936035156

# tx13 "Mixotroph" : (6333677,3602844) -> (6333721,3602575)
$sh = $g.GroupItems.Item(11)
$sh.Left = 498.71820068359375
$sh.Top  = 283.6673278808594

# tx14 " 13.115 %" : (6351568,3949773) -> (6351612,3949505)
$sh = $g.GroupItems.Item(12)
$sh.Left = 500.126953125
$sh.Top  = 310.9846496582031

# tx15 "Heterotroph" : (5026331,5397598) -> (5613741,4810124)
$sh = $g.GroupItems.Item(13)
$sh.Left = 442.02685546875
$sh.Top  = 378.74993896484375

# tx16 " 42.623 %" : (5140718,5744528) -> (5728127,5157054)
$sh = $g.GroupItems.Item(14)
$sh.Left = 451.03363037109375
$sh.Top  = 406.0672607421875

# tx17 "Autotrophic" : (3523550,3657693) -> (3523518,3657674)
$sh = $g.GroupItems.Item(15)
$sh.Left = 277.4423828125
$sh.Top  = 288.0058288574219

# tx18 " 34.426 %" : (3613853,4004623) -> (3613820,4004604)
$sh = $g.GroupItems.Item(16)
$sh.Left = 284.5527648925781
$sh.Top  = 315.3231506347656

